# This script updates the vm_pu (per-unit voltage magnitude) results sheet
# to reflect a fix in the Q (reactive power) calculation, per commit:
# "fixed issue anthony about Q calculation"
#
# It rewrites the numeric results in columns C-H and L-O for rows 2-25,
# and populates newly-added result columns I and J (which were previously
# absent) for the same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{R=2; C=3; V=[double]0.9842905732752075},
    @{R=2; C=4; V=[double]0.9811283385523365},
    @{R=2; C=5; V=[double]0.9698376372579449},
    @{R=2; C=7; V=[double]0.9512246635179655},
    @{R=2; C=8; V=[double]0.9820182892046919},
    @{R=2; C=9; V=[double]1},
    @{R=2; C=10; V=[double]0.9999999999999999},
    @{R=2; C=12; V=[double]0.9764981810933594},
    @{R=2; C=13; V=[double]0.9612184336118453},
    @{R=2; C=14; V=[double]0.9584380989516799},
    @{R=2; C=15; V=[double]0.943132282360777},
    @{R=3; C=3; V=[double]0.9847614301686957},
    @{R=3; C=4; V=[double]0.9816092177497895},
    @{R=3; C=5; V=[double]0.9707693143861075},
    @{R=3; C=6; V=[double]1},
    @{R=3; C=7; V=[double]0.9529905537392948},
    @{R=3; C=8; V=[double]0.9824996045925378},
    @{R=3; C=9; V=[double]0.9999999999999999},
    @{R=3; C=10; V=[double]1},
    @{R=3; C=12; V=[double]0.9770603059090565},
    @{R=3; C=13; V=[double]0.9617999667312049},
    @{R=3; C=14; V=[double]0.9598753472600494},
    @{R=3; C=15; V=[double]0.9450083754515111},
    @{R=4; C=3; V=[double]0.9850031324577605},
    @{R=4; C=4; V=[double]0.9819638406804874},
    @{R=4; C=5; V=[double]0.9715399486368698},
    @{R=4; C=6; V=[double]0.9999999999999998},
    @{R=4; C=7; V=[double]0.9543414143742412},
    @{R=4; C=8; V=[double]0.9828545491905407},
    @{R=4; C=9; V=[double]1},
    @{R=4; C=10; V=[double]0.9999999999999998},
    @{R=4; C=12; V=[double]0.9773715164757277},
    @{R=4; C=13; V=[double]0.9622319600074517},
    @{R=4; C=14; V=[double]0.9610302211027776},
    @{R=4; C=15; V=[double]0.9464437657231575},
    @{R=5; C=3; V=[double]0.9851270066711019},
    @{R=5; C=4; V=[double]0.9821867405528472},
    @{R=5; C=5; V=[double]0.9720502289313449},
    @{R=5; C=6; V=[double]1},
    @{R=5; C=7; V=[double]0.9550840420252859},
    @{R=5; C=8; V=[double]0.983077651248363},
    @{R=5; C=9; V=[double]1},
    @{R=5; C=10; V=[double]0.9999999999999999},
    @{R=5; C=12; V=[double]0.977531440146133},
    @{R=5; C=13; V=[double]0.9624958020644644},
    @{R=5; C=14; V=[double]0.9617387384972735},
    @{R=5; C=15; V=[double]0.9472304333549333},
    @{R=6; C=3; V=[double]0.9851920106479626},
    @{R=6; C=4; V=[double]0.9823025138319582},
    @{R=6; C=5; V=[double]0.9723160284762229},
    @{R=6; C=7; V=[double]0.9554250139826906},
    @{R=6; C=8; V=[double]0.9831935295417734},
    @{R=6; C=9; V=[double]0.9999999999999998},
    @{R=6; C=10; V=[double]0.9999999999999998},
    @{R=6; C=12; V=[double]0.9776120882110791},
    @{R=6; C=13; V=[double]0.9626296441108028},
    @{R=6; C=14; V=[double]0.9620900712332151},
    @{R=6; C=15; V=[double]0.9475905412893052},
    @{R=7; C=3; V=[double]0.9851920106479626},
    @{R=7; C=4; V=[double]0.9823025138319582},
    @{R=7; C=5; V=[double]0.9723160284762229},
    @{R=7; C=7; V=[double]0.9554250139826906},
    @{R=7; C=8; V=[double]0.9831935295417734},
    @{R=7; C=9; V=[double]0.9999999999999998},
    @{R=7; C=10; V=[double]0.9999999999999998},
    @{R=7; C=12; V=[double]0.9776120882110791},
    @{R=7; C=13; V=[double]0.9626296441108028},
    @{R=7; C=14; V=[double]0.9620900712332151},
    @{R=7; C=15; V=[double]0.9475905412893052},
    @{R=8; C=3; V=[double]0.9851963447129699},
    @{R=8; C=4; V=[double]0.9822976981034642},
    @{R=8; C=5; V=[double]0.9723008864926829},
    @{R=8; C=6; V=[double]1},
    @{R=8; C=7; V=[double]0.9553654879880982},
    @{R=8; C=8; V=[double]0.9831887094450836},
    @{R=8; C=9; V=[double]0.9999999999999999},
    @{R=8; C=10; V=[double]1},
    @{R=8; C=12; V=[double]0.9776127300166201},
    @{R=8; C=13; V=[double]0.9626208767159098},
    @{R=8; C=14; V=[double]0.962054607473355},
    @{R=8; C=15; V=[double]0.9475265134078941},
    @{R=9; C=3; V=[double]0.9850903731009196},
    @{R=9; C=4; V=[double]0.9820719060691018},
    @{R=9; C=5; V=[double]0.9717786797404243},
    @{R=9; C=6; V=[double]0.9999999999999999},
    @{R=9; C=7; V=[double]0.954336133454451},
    @{R=9; C=8; V=[double]0.9829627126018696},
    @{R=9; C=9; V=[double]1},
    @{R=9; C=10; V=[double]1},
    @{R=9; C=12; V=[double]0.9774514275959572},
    @{R=9; C=13; V=[double]0.9623340173835342},
    @{R=9; C=14; V=[double]0.9612285952560575},
    @{R=9; C=15; V=[double]0.9464298190309282},
    @{R=10; C=3; V=[double]0.9843638546829151},
    @{R=10; C=4; V=[double]0.9812498842331583},
    @{R=10; C=5; V=[double]0.9701197462836108},
    @{R=10; C=6; V=[double]0.9999999999999998},
    @{R=10; C=7; V=[double]0.9514460344886316},
    @{R=10; C=8; V=[double]0.9821399451357771},
    @{R=10; C=9; V=[double]1},
    @{R=10; C=10; V=[double]0.9999999999999999},
    @{R=10; C=12; V=[double]0.9765786327527867},
    @{R=10; C=13; V=[double]0.9613493920483031},
    @{R=10; C=14; V=[double]0.9587622369341264},
    @{R=10; C=15; V=[double]0.9433627764895385},
    @{R=11; C=3; V=[double]0.983150898026499},
    @{R=11; C=4; V=[double]0.9802397245915735},
    @{R=11; C=5; V=[double]0.9683209464676763},
    @{R=11; C=6; V=[double]1},
    @{R=11; C=7; V=[double]0.9483517822023083},
    @{R=11; C=8; V=[double]0.9811288692101582},
    @{R=11; C=9; V=[double]0.9999999999999999},
    @{R=11; C=10; V=[double]0.9999999999999999},
    @{R=11; C=12; V=[double]0.9752137560304796},
    @{R=11; C=13; V=[double]0.9601709278226146},
    @{R=11; C=14; V=[double]0.9561195670059685},
    @{R=11; C=15; V=[double]0.9400862169863465},
    @{R=12; C=3; V=[double]0.9824172345915887},
    @{R=12; C=4; V=[double]0.9797054575049613},
    @{R=12; C=5; V=[double]0.967446583497253},
    @{R=12; C=6; V=[double]0.9999999999999998},
    @{R=12; C=7; V=[double]0.9467899189745098},
    @{R=12; C=8; V=[double]0.9805941175066785},
    @{R=12; C=9; V=[double]1},
    @{R=12; C=10; V=[double]0.9999999999999997},
    @{R=12; C=12; V=[double]0.9744038021747637},
    @{R=12; C=13; V=[double]0.9595527172830182},
    @{R=12; C=14; V=[double]0.9548200488617202},
    @{R=12; C=15; V=[double]0.9384332650196805},
    @{R=13; C=3; V=[double]0.9820023940045377},
    @{R=13; C=4; V=[double]0.9793299024301726},
    @{R=13; C=5; V=[double]0.966786926731754},
    @{R=13; C=7; V=[double]0.9456037364257021},
    @{R=13; C=8; V=[double]0.9802182217776954},
    @{R=13; C=9; V=[double]0.9999999999999999},
    @{R=13; C=10; V=[double]0.9999999999999999},
    @{R=13; C=12; V=[double]0.9739300288623493},
    @{R=13; C=13; V=[double]0.9591117086934322},
    @{R=13; C=14; V=[double]0.9538204694223799},
    @{R=13; C=15; V=[double]0.9371750876344642},
    @{R=14; C=3; V=[double]0.9820878612303036},
    @{R=14; C=4; V=[double]0.9793223922251755},
    @{R=14; C=5; V=[double]0.9667386462022428},
    @{R=14; C=7; V=[double]0.9454559575221704},
    @{R=14; C=8; V=[double]0.9802107047604277},
    @{R=14; C=9; V=[double]0.9999999999999999},
    @{R=14; C=10; V=[double]0.9999999999999999},
    @{R=14; C=12; V=[double]0.9740065569032638},
    @{R=14; C=13; V=[double]0.9590939472817183},
    @{R=14; C=14; V=[double]0.9537132446086847},
    @{R=14; C=15; V=[double]0.9370151759200641},
    @{R=15; C=3; V=[double]0.9823056960329208},
    @{R=15; C=4; V=[double]0.9794752703124122},
    @{R=15; C=5; V=[double]0.967002473493154},
    @{R=15; C=6; V=[double]1},
    @{R=15; C=7; V=[double]0.9458182784858243},
    @{R=15; C=8; V=[double]0.980363721518568},
    @{R=15; C=9; V=[double]0.9999999999999999},
    @{R=15; C=10; V=[double]0.9999999999999999},
    @{R=15; C=12; V=[double]0.9742404544516418},
    @{R=15; C=13; V=[double]0.9592647976845325},
    @{R=15; C=14; V=[double]0.9540662897119947},
    @{R=15; C=15; V=[double]0.9373965007352465},
    @{R=16; C=3; V=[double]0.9823374137063791},
    @{R=16; C=4; V=[double]0.979461474895316},
    @{R=16; C=5; V=[double]0.966941949677253},
    @{R=16; C=6; V=[double]1},
    @{R=16; C=7; V=[double]0.9458018467489353},
    @{R=16; C=8; V=[double]0.9803499135880829},
    @{R=16; C=9; V=[double]0.9999999999999997},
    @{R=16; C=10; V=[double]0.9999999999999999},
    @{R=16; C=12; V=[double]0.9742723345537855},
    @{R=16; C=13; V=[double]0.9592505526722926},
    @{R=16; C=14; V=[double]0.9540036222581363},
    @{R=16; C=15; V=[double]0.9373797552456503},
    @{R=17; C=3; V=[double]0.9823555207486688},
    @{R=17; C=4; V=[double]0.9794335327439438},
    @{R=17; C=5; V=[double]0.9668510594601664},
    @{R=17; C=7; V=[double]0.945760194684361},
    @{R=17; C=8; V=[double]0.9803219460912642},
    @{R=17; C=9; V=[double]0.9999999999999998},
    @{R=17; C=10; V=[double]0.9999999999999999},
    @{R=17; C=12; V=[double]0.9742897575009037},
    @{R=17; C=13; V=[double]0.9592210966808317},
    @{R=17; C=14; V=[double]0.9539056389477143},
    @{R=17; C=15; V=[double]0.937336747073798},
    @{R=18; C=3; V=[double]0.9827867998449417},
    @{R=18; C=4; V=[double]0.9797765586449896},
    @{R=18; C=5; V=[double]0.9674405499474078},
    @{R=18; C=7; V=[double]0.9467684789817361},
    @{R=18; C=8; V=[double]0.9806652831403141},
    @{R=18; C=9; V=[double]1},
    @{R=18; C=10; V=[double]0.9999999999999999},
    @{R=18; C=12; V=[double]0.9747700908107946},
    @{R=18; C=13; V=[double]0.9596184939631849},
    @{R=18; C=14; V=[double]0.9547745343684303},
    @{R=18; C=15; V=[double]0.9384043402627187},
    @{R=19; C=3; V=[double]0.9830032890252697},
    @{R=19; C=4; V=[double]0.9799577675828783},
    @{R=19; C=5; V=[double]0.9677475123646305},
    @{R=19; C=6; V=[double]1},
    @{R=19; C=7; V=[double]0.9473556670944514},
    @{R=19; C=8; V=[double]0.98084665644713},
    @{R=19; C=9; V=[double]1},
    @{R=19; C=10; V=[double]1},
    @{R=19; C=12; V=[double]0.975016156709174},
    @{R=19; C=13; V=[double]0.9598323277670746},
    @{R=19; C=14; V=[double]0.9552502288399704},
    @{R=19; C=15; V=[double]0.9390274224138258},
    @{R=20; C=3; V=[double]0.9828363654067928},
    @{R=20; C=4; V=[double]0.9799168850626159},
    @{R=20; C=5; V=[double]0.9677252194202214},
    @{R=20; C=6; V=[double]0.9999999999999998},
    @{R=20; C=7; V=[double]0.9473098458567025},
    @{R=20; C=8; V=[double]0.9808057368436198},
    @{R=20; C=9; V=[double]0.9999999999999999},
    @{R=20; C=10; V=[double]1},
    @{R=20; C=12; V=[double]0.9748475017019447},
    @{R=20; C=13; V=[double]0.9597903213236832},
    @{R=20; C=14; V=[double]0.9552262777431489},
    @{R=20; C=15; V=[double]0.9389809226738215},
    @{R=21; C=3; V=[double]0.9817879768660259},
    @{R=21; C=4; V=[double]0.9792444082730241},
    @{R=21; C=5; V=[double]0.9666764365630872},
    @{R=21; C=6; V=[double]0.9999999999999999},
    @{R=21; C=7; V=[double]0.9455600584940287},
    @{R=21; C=8; V=[double]0.980132650071486},
    @{R=21; C=9; V=[double]0.9999999999999999},
    @{R=21; C=10; V=[double]1},
    @{R=21; C=12; V=[double]0.9737168284953315},
    @{R=21; C=13; V=[double]0.9590276185139186},
    @{R=21; C=14; V=[double]0.9537266232386542},
    @{R=21; C=15; V=[double]0.9371344723387238},
    @{R=22; C=3; V=[double]0.9808279762086934},
    @{R=22; C=4; V=[double]0.9785705444238441},
    @{R=22; C=5; V=[double]0.9655805507542865},
    @{R=22; C=6; V=[double]1},
    @{R=22; C=7; V=[double]0.9438134485673872},
    @{R=22; C=8; V=[double]0.9794581749816066},
    @{R=22; C=9; V=[double]0.9999999999999998},
    @{R=22; C=10; V=[double]0.9999999999999998},
    @{R=22; C=12; V=[double]0.9726754884541766},
    @{R=22; C=13; V=[double]0.9582637016743164},
    @{R=22; C=14; V=[double]0.9521744115393942},
    @{R=22; C=15; V=[double]0.9352909368162606},
    @{R=23; C=3; V=[double]0.9800518008275102},
    @{R=23; C=4; V=[double]0.9780215174746277},
    @{R=23; C=5; V=[double]0.9646824591353366},
    @{R=23; C=7; V=[double]0.9424238660135128},
    @{R=23; C=8; V=[double]0.9789086500273156},
    @{R=23; C=9; V=[double]0.9999999999999998},
    @{R=23; C=10; V=[double]0.9999999999999999},
    @{R=23; C=12; V=[double]0.9718353181777779},
    @{R=23; C=13; V=[double]0.9576436476299582},
    @{R=23; C=14; V=[double]0.9509141089728387},
    @{R=23; C=15; V=[double]0.9338248413287749},
    @{R=24; C=3; V=[double]0.97976407995092},
    @{R=24; C=4; V=[double]0.9778043429007098},
    @{R=24; C=5; V=[double]0.964309909663876},
    @{R=24; C=7; V=[double]0.9418894930470122},
    @{R=24; C=8; V=[double]0.9786912784611733},
    @{R=24; C=9; V=[double]0.9999999999999999},
    @{R=24; C=10; V=[double]0.9999999999999998},
    @{R=24; C=12; V=[double]0.9715235152983457},
    @{R=24; C=13; V=[double]0.9573994685050894},
    @{R=24; C=14; V=[double]0.9504033447476451},
    @{R=24; C=15; V=[double]0.9332614349251654},
    @{R=25; C=3; V=[double]0.9813346941439727},
    @{R=25; C=4; V=[double]0.9788033680055725},
    @{R=25; C=5; V=[double]0.9658069236334917},
    @{R=25; C=6; V=[double]0.9999999999999999},
    @{R=25; C=7; V=[double]0.9445172624681694},
    @{R=25; C=8; V=[double]0.9796912097502821},
    @{R=25; C=9; V=[double]1},
    @{R=25; C=10; V=[double]1},
    @{R=25; C=12; V=[double]0.9732189799086438},
    @{R=25; C=13; V=[double]0.9585348889415},
    @{R=25; C=14; V=[double]0.9525974055419647},
    @{R=25; C=15; V=[double]0.9360369856595493}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.R, $u.C).Value = $u.V
}

Write-Host ("Updated {0} cells" -f $updates.Count)
